$d = $word.ActiveDocument
$vtab = [char]11

# ---------------------------------------------------------------------------
# 1. Title paragraph: turn
#      "Encryption (JWE) C# Library<br>.NET 4.8 Solution<br>.NET Standard 2.0 / Core 3.1 Solution<br>"
#    into
#      "JSON Web Encryption (JWE) using<br>C# and .NET 4.8 Solution or<br>.NET Standard 2.0 / Core 3.1 Solution<br>"
# ---------------------------------------------------------------------------

# "Encryption (JWE)" -> "JSON Web Encryption (JWE)"
$d.Content.Find.Execute("Encryption (JWE)", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "JSON Web Encryption (JWE)", 2) | Out-Null

# " C#" -> " using<br>C#"
$d.Content.Find.Execute(" C#", $true, $false, $false, $false, $false, `
                         $true, 1, $false, (" using" + $vtab + "C#"), 2) | Out-Null

# "Library<br>" -> "and "
$d.Content.Find.Execute(("Library" + $vtab), $true, $false, $false, $false, $false, `
                         $true, 1, $false, "and ", 2) | Out-Null

# ".NET 4.8 Solution<br>.NET Standard" -> ".NET 4.8 Solution or<br>.NET Standard"
$d.Content.Find.Execute((".NET 4.8 Solution" + $vtab + ".NET Standard"), $true, $false, $false, $false, $false, `
                         $true, 1, $false, (".NET 4.8 Solution or" + $vtab + ".NET Standard"), 2) | Out-Null

# Insert a (collapsed) "_GoBack" bookmark right after the newly-inserted " or".
# Word bookmark names are unique, so re-adding "_GoBack" here automatically
# relocates the existing one (previously at the end of the document) to this
# new position - matching the diff, which both removes it from the last
# paragraph and adds it here.
$findRange = $d.Content
$findRange.Find.Execute(" or", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "", 0) | Out-Null
$bmRange = $d.Range($findRange.End, $findRange.End)
$d.Bookmarks.Add("_GoBack", $bmRange) | Out-Null

# ---------------------------------------------------------------------------
# 2. Final bullet: "Uses the same Key Id to decrypt string using private key"
#    -> "Uses public Key Id is the same as the private Key Id."
# ---------------------------------------------------------------------------
$d.Content.Find.Execute("Uses the same Key Id to decrypt string using private key", `
                         $true, $false, $false, $false, $false, `
                         $true, 1, $false, "Uses public Key Id is the same as the private Key Id.", 2) | Out-Null
